# "sorting output for strategytree"
# Adds a second, alternately-formatted copy of the CardNumberOrScore lookup
# table to Sheet3 (rows 27:60), plus a stray scratch value in J25, and
# updates the saved selections on Sheet1/Sheet3.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# Row 25: a lone helper/scratch cell in column J holding the literal text
# that the first new formula (G27) below will also compute.
$ws3.Cells.Item(25, 10).Value2 = "CardNumberOrScore.Ace => 11,"

# New table: same 34 labels as B4:D24 (Ace..King, Score21..Score2, Busted)
# but rebuilt further down the sheet (rows 27-60) using column E for the
# numeric operand and a slightly different separator in the formula
# ("label => value," with no literal quotes around the number).
$rows = @(
  @(27, "Ace",     11),
  @(28, "Two",      2),
  @(29, "Three",    3),
  @(30, "Four",     4),
  @(31, "Five",     5),
  @(32, "Six",      6),
  @(33, "Seven",    7),
  @(34, "Eight",    8),
  @(35, "Nine",     9),
  @(36, "Ten",     10),
  @(37, "Jack",    10),
  @(38, "Queen",   10),
  @(39, "King",    10),
  @(40, "Score21", 21),
  @(41, "Score20", 20),
  @(42, "Score19", 19),
  @(43, "Score18", 18),
  @(44, "Score17", 17),
  @(45, "Score16", 16),
  @(46, "Score15", 15),
  @(47, "Score14", 14),
  @(48, "Score13", 13),
  @(49, "Score12", 12),
  @(50, "Score11", 11),
  @(51, "Score10", 10),
  @(52, "Score9",   9),
  @(53, "Score8",   8),
  @(54, "Score7",   7),
  @(55, "Score6",   6),
  @(56, "Score5",   5),
  @(57, "Score4",   4),
  @(58, "Score3",   3),
  @(59, "Score2",   2),
  @(60, "Busted",  22)
)

foreach ($item in $rows) {
  $r = $item[0]
  $ws3.Cells.Item($r, 2).Value2 = "CardNumberOrScore."
  $ws3.Cells.Item($r, 3).Value2 = $item[1]
  $ws3.Cells.Item($r, 5).Value2 = $item[2]
}

# G27 is its own (non-shared) formula; G28:G60 share one formula group.
$ws3.Range("G27").Formula = '=CONCAT(B27,C27," => ",E27,",")'
$ws3.Range("G28:G60").Formula = '=CONCAT(B28,C28," => ",E28,",")'

# Restore the sheet selections to match what was saved (Sheet1's selection
# moves to B15; Sheet3 ends up re-selected over the new formula column).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B15").Select()

$ws3.Activate()
$ws3.Range("G27:G60").Select()
